$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.293.08"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "2.097.50"
$ws.Range("E3").Value = "  +4.60%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +21.93%  "
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.375"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0744"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.45%  "
$ws.Range("E12").Value = "  +8.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.61%  "
$ws.Range("D14").Value = "2.399.88"
$ws.Range("E14").Value = "  +4.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.841"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.12%  "
$ws.Range("D16").Value = "2.092.66"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.00%  "
$ws.Range("D18").Value = "37.295.75"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +15.33%  "
$ws.Range("D21").Value = "0.0₃0845"
$ws.Range("E21").Value = "  +3.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.92%  "
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("B31").Value = "Gas"
$ws.Range("C31").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +27.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0618"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0906"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.95%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.86%  "
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0226"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0928"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +97.32%  "
$ws.Range("D48").Value = "1.321.39"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("E49").Value = "  +6.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.63%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.48%  "